$wb = $excel.ActiveWorkbook

$pirData = @(
    @('2026-02-06','09:52:39','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:52:40','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:52:44','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:52:49','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:52:54','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:52:59','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:53:04','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:53:09','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:53:14','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:53:19','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:53:24','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:53:29','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:53:30','09:00','Bathroom','Motion Detected','Active'),
    @('2026-02-06','09:53:38','09:00','Bathroom','No Motion','Inactive')
)

$humidityData = @(
    @('2026-02-06','09:52:39','09:00','Bathroom','71.5%','Active'),
    @('2026-02-06','09:52:41','09:00','Bathroom','71.4%','Active'),
    @('2026-02-06','09:52:46','09:00','Bathroom','71.4%','Active'),
    @('2026-02-06','09:52:52','09:00','Bathroom','71.4%','Active'),
    @('2026-02-06','09:52:56','09:00','Bathroom','71.5%','Active'),
    @('2026-02-06','09:53:01','09:00','Bathroom','71.3%','Active'),
    @('2026-02-06','09:53:07','09:00','Bathroom','70.2%','Active'),
    @('2026-02-06','09:53:12','09:00','Bathroom','71.2%','Active'),
    @('2026-02-06','09:53:17','09:00','Bathroom','70.2%','Active'),
    @('2026-02-06','09:53:27','09:00','Bathroom','70.1%','Active'),
    @('2026-02-06','09:53:31','09:00','Bathroom','71.0%','Active'),
    @('2026-02-06','09:53:37','09:00','Bathroom','70.2%','Active')
)

$temperatureData = @(
    @('2026-02-06','09:52:39','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:52:42','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:52:47','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:52:52','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:52:57','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:53:02','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:53:08','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:53:12','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:53:17','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:53:27','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:53:32','09:00','Bathroom','27.7C','Active'),
    @('2026-02-06','09:53:37','09:00','Bathroom','27.8C','Active')
)

$wsPIR = $wb.Worksheets.Item("PIR")
$startRowPIR = 151
$wsPIR.Range("A$startRowPIR`:A$($startRowPIR + $pirData.Count - 1)").NumberFormat = "@"
$r = $startRowPIR
foreach ($row in $pirData) {
    $wsPIR.Range("A$r").Value = $row[0]
    $wsPIR.Range("B$r").Value = $row[1]
    $wsPIR.Range("C$r").Value = $row[2]
    $wsPIR.Range("D$r").Value = $row[3]
    $wsPIR.Range("E$r").Value = $row[4]
    $wsPIR.Range("F$r").Value = $row[5]
    $r = $r + 1
}

$wsHumidity = $wb.Worksheets.Item("Humidity")
$startRowHumidity = 80
$wsHumidity.Range("A$startRowHumidity`:A$($startRowHumidity + $humidityData.Count - 1)").NumberFormat = "@"
$wsHumidity.Range("E$startRowHumidity`:E$($startRowHumidity + $humidityData.Count - 1)").NumberFormat = "@"
$r = $startRowHumidity
foreach ($row in $humidityData) {
    $wsHumidity.Range("A$r").Value = $row[0]
    $wsHumidity.Range("B$r").Value = $row[1]
    $wsHumidity.Range("C$r").Value = $row[2]
    $wsHumidity.Range("D$r").Value = $row[3]
    $wsHumidity.Range("E$r").Value = $row[4]
    $wsHumidity.Range("F$r").Value = $row[5]
    $r = $r + 1
}

$wsTemperature = $wb.Worksheets.Item("Temperature")
$startRowTemperature = 80
$wsTemperature.Range("A$startRowTemperature`:A$($startRowTemperature + $temperatureData.Count - 1)").NumberFormat = "@"
$r = $startRowTemperature
foreach ($row in $temperatureData) {
    $wsTemperature.Range("A$r").Value = $row[0]
    $wsTemperature.Range("B$r").Value = $row[1]
    $wsTemperature.Range("C$r").Value = $row[2]
    $wsTemperature.Range("D$r").Value = $row[3]
    $wsTemperature.Range("E$r").Value = $row[4]
    $wsTemperature.Range("F$r").Value = $row[5]
    $r = $r + 1
}
